$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02950307763024241
$ws.Range("D2").Value = 0.1292627146720763
$ws.Range("G2").Value = 0.1248244242667473
$ws.Range("H2").Value = 0.99
